# Applies the diff:
#  1. Updates the statistics for row 25 ("3/4 Salto vw A") with recalculated values.
#  2. Inserts a new row before old row 34, shifting rows 34-51 down to 35-52,
#     and populates the new row 34 with data for "3/4 Salto Vw A".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: update row 25 values (C25:X25) ---
$ws.Range("C25").Value = 14.811
$ws.Range("D25").Value = 14.7015
$ws.Range("E25").Value = -1.875375
$ws.Range("F25").Value = -0.6473749999999999
$ws.Range("G25").Value = 5.05025
$ws.Range("H25").Value = 8.1205
$ws.Range("I25").Value = -36.26512499999999
$ws.Range("J25").Value = 13.29325
$ws.Range("K25").Value = 13.04425
$ws.Range("L25").Value = -66.11762499999999
$ws.Range("M25").Value = 20.946125
$ws.Range("N25").Value = 34.19675
$ws.Range("O25").Value = 304.482625
$ws.Range("P25").Value = 53.268125
$ws.Range("Q25").Value = 27.243875
$ws.Range("R25").Value = 23.109625
$ws.Range("S25").Value = 16.093875
$ws.Range("T25").Value = 1.971625
$ws.Range("U25").Value = 20.625375
$ws.Range("V25").Value = 13.86125
$ws.Range("W25").Value = 123.710125
$ws.Range("X25").Value = 23.838875

# --- Step 2: insert a new row before row 34 ---
$ws.Rows.Item(34).Insert()

# Restore the formatting of column A on the newly inserted row (matches the
# style used for column A on every other data row: bold, centered, bordered).
$ws.Cells.Item(34, 1).Value = 0
$ws.Cells.Item(34, 1).Font.Bold = $true
$ws.Cells.Item(34, 1).HorizontalAlignment = -4108
$ws.Cells.Item(34, 1).VerticalAlignment = -4160
$ws.Cells.Item(34, 1).Borders.LineStyle = 1

# --- Step 3: populate the new row 34 with "3/4 Salto Vw A" data ---
$ws.Range("B34").Value = "3/4 Salto Vw A"
$ws.Range("C34").Value = 14.8887
$ws.Range("D34").Value = 14.6743
$ws.Range("E34").Value = -5.4937
$ws.Range("F34").Value = -0.4318
$ws.Range("G34").Value = 4.5641
$ws.Range("H34").Value = 2.2727
$ws.Range("I34").Value = -157.2026
$ws.Range("J34").Value = 10.04
$ws.Range("K34").Value = 4.5112
$ws.Range("L34").Value = -276.9653
$ws.Range("M34").Value = 15.114
$ws.Range("N34").Value = 22.2219
$ws.Range("O34").Value = 305.0548
$ws.Range("P34").Value = 40.2427
$ws.Range("Q34").Value = 27.9702
$ws.Range("R34").Value = 23.3319
$ws.Range("S34").Value = 17.549
$ws.Range("T34").Value = 1.383
$ws.Range("U34").Value = 19.9329
$ws.Range("V34").Value = 12.1583
$ws.Range("W34").Value = 121.5812
$ws.Range("X34").Value = 17.0754

Write-Host "Edit applied successfully"
